# Weekly update for Camote (Vega Central Mapocho de Santiago) sheet.
# Two new weekly observation rows are inserted at rows 54-55 (pushing the
# existing rows 54-94 down to rows 56-96), and populated with new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the top of the insertion point; existing rows
# 54-94 (and everything below) shift down to 56-96.
$ws.Rows("54:55").Insert()

# ---- New row 54 -------------------------------------------------------
$ws.Range("A54").Value = 9
$ws.Range("B54").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C54").Value = "Metropolitana"
$ws.Range("D54").Value = 44767
$ws.Range("E54").Value = 13
$ws.Range("F54").Value = 100114002
$ws.Range("G54").Value = "Camote"
$ws.Range("H54").Value = "Sin especificar"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 430
$ws.Range("K54").Value = 13000
$ws.Range("L54").Value = 14000
$ws.Range("M54").Value = 13500
$ws.Range("N54").Value = "$/caja 18 kilos"
$ws.Range("O54").Value = "Perú"
$ws.Range("P54").Value = 750
$ws.Range("Q54").Value = 18
$ws.Range("R54").Value = "Hortaliza"

# ---- New row 55 -------------------------------------------------------
$ws.Range("A55").Value = 9
$ws.Range("B55").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C55").Value = "Metropolitana"
$ws.Range("D55").Value = 44767
$ws.Range("E55").Value = 13
$ws.Range("F55").Value = 100114002
$ws.Range("G55").Value = "Camote"
$ws.Range("H55").Value = "Sin especificar"
$ws.Range("I55").Value = "Primera"
$ws.Range("J55").Value = 880
$ws.Range("K55").Value = 9000
$ws.Range("L55").Value = 10000
$ws.Range("M55").Value = 9500
$ws.Range("N55").Value = "$/malla 18 kilos"
$ws.Range("O55").Value = "Perú"
$ws.Range("P55").Value = 528
$ws.Range("Q55").Value = 18
$ws.Range("R55").Value = "Hortaliza"
